$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 587.5
$ws.Range("I9").Value = 483.33334
$ws.Range("J9").Value = 900
$ws.Range("K9").Value = 483.33334
$ws.Range("L9").Value = 900
$ws.Range("M9").Value = -314.33334
$ws.Range("N9").Value = -1238
$ws.Range("H11").Value = 407.38095
$ws.Range("I11").Value = 407.38095
$ws.Range("K11").Value = 407.38095
$ws.Range("M11").Value = -267.38095
$ws.Range("H32").Value = 4129.381
$ws.Range("I32").Value = 3893.818
$ws.Range("J32").Value = 4388.5
$ws.Range("K32").Value = 3893.818
$ws.Range("L32").Value = 4388.5
$ws.Range("M32").Value = -3567.818
$ws.Range("N32").Value = -5040.5
$ws.Range("H39").Value = 2073.2222
$ws.Range("I39").Value = 2560.125
$ws.Range("K39").Value = 7680.375
$ws.Range("M39").Value = -7384.375
$ws.Range("H132").Value = 3613.4644
$ws.Range("I132").Value = 1433.875
$ws.Range("K132").Value = 4301.625
$ws.Range("M132").Value = -1771.625
$ws.Range("H135").Value = 2155.5
$ws.Range("I135").Value = 2155.5
$ws.Range("K135").Value = 19399.5
$ws.Range("M135").Value = -16864.5
$ws.Range("H138").Value = 3239.9363
$ws.Range("I138").Value = 2184.4546
$ws.Range("J138").Value = 3562.4443
$ws.Range("K138").Value = 6553.3638
$ws.Range("L138").Value = 10687.3329
$ws.Range("M138").Value = -1413.3638
$ws.Range("N138").Value = -20967.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 1000.3333
$ws.Range("J6").Value = 1000.3333
$ws.Range("L6").Value = 1000.3333
$ws.Range("N6").Value = -1346.3333
$ws.Range("H61").Value = 1802.1143
$ws.Range("I61").Value = 1747.129
$ws.Range("K61").Value = 1747.129
$ws.Range("M61").Value = -1535.129
$ws.Range("H63").Value = 4278.3
$ws.Range("I63").Value = 2505.077
$ws.Range("K63").Value = 2505.077
$ws.Range("M63").Value = -1819.077
$ws.Range("H66").Value = 4278.3
$ws.Range("I66").Value = 2505.077
$ws.Range("K66").Value = 12525.385
$ws.Range("M66").Value = -9093.385000000002
$ws.Range("H96").Value = 30000
$ws.Range("J96").Value = 30000
$ws.Range("L96").Value = 30000
$ws.Range("N96").Value = -35492
$ws.Range("H122").Value = 2096.5945
$ws.Range("I122").Value = 1545.2963
$ws.Range("K122").Value = 4635.8889
$ws.Range("M122").Value = -2185.8889
$ws.Range("H132").Value = 2325.9375
$ws.Range("I132").Value = 1942.9756
$ws.Range("K132").Value = 5828.9268
$ws.Range("M132").Value = -3298.9268
$ws.Range("H136").Value = 1802.1143
$ws.Range("I136").Value = 1747.129
$ws.Range("K136").Value = 5241.387
$ws.Range("M136").Value = -2691.387

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 25007632
$ws.Range("I20").Value = 27785668
$ws.Range("K20").Value = 27785668
$ws.Range("M20").Value = -27785421
$ws.Range("H80").Value = 705.2
$ws.Range("I80").Value = 498.6
$ws.Range("J80").Value = 808.5
$ws.Range("K80").Value = 498.6
$ws.Range("L80").Value = 808.5
$ws.Range("M80").Value = 499.4
$ws.Range("N80").Value = -2804.5
$ws.Range("H83").Value = 705.2
$ws.Range("I83").Value = 498.6
$ws.Range("J83").Value = 808.5
$ws.Range("K83").Value = 2493
$ws.Range("L83").Value = 4042.5
$ws.Range("M83").Value = 2499
$ws.Range("N83").Value = -14026.5
$ws.Range("H95").Value = 60416
$ws.Range("J95").Value = 60416
$ws.Range("L95").Value = 60416
$ws.Range("N95").Value = -65908
$ws.Range("H99").Value = 2498.8635
$ws.Range("I99").Value = 2587.6316
$ws.Range("K99").Value = 2587.6316
$ws.Range("M99").Value = -1089.6316

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31413.285
$ws.Range("I31").Value = 38453.85
$ws.Range("J31").Value = 7651.375
$ws.Range("K31").Value = 38453.85
$ws.Range("L31").Value = 7651.375
$ws.Range("M31").Value = -38158.85
$ws.Range("N31").Value = -8241.375
$ws.Range("H33").Value = 1150
$ws.Range("I33").Value = 1150
$ws.Range("K33").Value = 1150
$ws.Range("M33").Value = -771
$ws.Range("H34").Value = 31413.285
$ws.Range("I34").Value = 38453.85
$ws.Range("J34").Value = 7651.375
$ws.Range("K34").Value = 38453.85
$ws.Range("L34").Value = 7651.375
$ws.Range("M34").Value = -38251.85
$ws.Range("N34").Value = -8055.375
$ws.Range("H35").Value = 4570.6924
$ws.Range("I35").Value = 2333
$ws.Range("J35").Value = 5242
$ws.Range("K35").Value = 2333
$ws.Range("L35").Value = 5242
$ws.Range("M35").Value = -2039
$ws.Range("N35").Value = -5830
$ws.Range("H86").Value = 3095.389
$ws.Range("I86").Value = 2848.0833
$ws.Range("J86").Value = 3590
$ws.Range("K86").Value = 2848.0833
$ws.Range("L86").Value = 3590
$ws.Range("M86").Value = -1725.0833
$ws.Range("N86").Value = -5836
$ws.Range("H89").Value = 3095.389
$ws.Range("I89").Value = 2848.0833
$ws.Range("J89").Value = 3590
$ws.Range("K89").Value = 14240.4165
$ws.Range("L89").Value = 17950
$ws.Range("M89").Value = -8624.416499999999
$ws.Range("N89").Value = -29182
$ws.Range("H99").Value = 15987.546
$ws.Range("I99").Value = 24301.5
$ws.Range("J99").Value = 6010.8
$ws.Range("K99").Value = 24301.5
$ws.Range("L99").Value = 6010.8
$ws.Range("M99").Value = -22803.5
$ws.Range("N99").Value = -9006.799999999999
$ws.Range("H105").Value = 1606.138
$ws.Range("I105").Value = 1440.7916
$ws.Range("J105").Value = 2399.8
$ws.Range("K105").Value = 1440.7916
$ws.Range("L105").Value = 2399.8
$ws.Range("M105").Value = 306.2084
$ws.Range("N105").Value = -5893.8
$ws.Range("H126").Value = 15987.546
$ws.Range("I126").Value = 24301.5
$ws.Range("J126").Value = 6010.8
$ws.Range("K126").Value = 72904.5
$ws.Range("L126").Value = 18032.4
$ws.Range("M126").Value = -70434.5
$ws.Range("N126").Value = -22972.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1411.5264
$ws.Range("I5").Value = 476
$ws.Range("K5").Value = 1428
$ws.Range("M5").Value = -1316
$ws.Range("H68").Value = 1263.5883
$ws.Range("I68").Value = 1011.5714
$ws.Range("K68").Value = 3034.7142
$ws.Range("M68").Value = -2223.7142
$ws.Range("H71").Value = 1263.5883
$ws.Range("I71").Value = 1011.5714
$ws.Range("K71").Value = 9104.142600000001
$ws.Range("M71").Value = -5048.142600000001
$ws.Range("H107").Value = 785.17645
$ws.Range("J107").Value = 789.9666999999999
$ws.Range("L107").Value = 2369.9001
$ws.Range("N107").Value = -6209.9001
$ws.Range("H129").Value = 71528.07000000001
$ws.Range("I129").Value = 334370.16
$ws.Range("J129").Value = 2960.5652
$ws.Range("K129").Value = 1003110.48
$ws.Range("L129").Value = 8881.695599999999
$ws.Range("M129").Value = -998110.48
$ws.Range("N129").Value = -18881.6956
$ws.Range("H131").Value = 1977.8292
$ws.Range("J131").Value = 1975.1714
$ws.Range("L131").Value = 5925.5142
$ws.Range("N131").Value = -16005.5142
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 27000
$ws.Range("N132").Value = -32060
$ws.Range("M132").ClearContents()
$ws.Range("H135").Value = 1411.5264
$ws.Range("I135").Value = 476
$ws.Range("K135").Value = 4284
$ws.Range("M135").Value = -1749

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7680.273
$ws.Range("I70").Value = 7685.5
$ws.Range("J70").Value = 7666.3335
$ws.Range("K70").Value = 7685.5
$ws.Range("L70").Value = 7666.3335
$ws.Range("M70").Value = -7415.5
$ws.Range("N70").Value = -8206.333500000001
$ws.Range("H73").Value = 7680.273
$ws.Range("I73").Value = 7685.5
$ws.Range("J73").Value = 7666.3335
$ws.Range("K73").Value = 7685.5
$ws.Range("L73").Value = 7666.3335
$ws.Range("M73").Value = -6749.5
$ws.Range("N73").Value = -9538.333500000001
$ws.Range("H113").Value = 2631.4443
$ws.Range("I113").Value = 2579.1667
$ws.Range("J113").Value = 2736
$ws.Range("K113").Value = 2579.1667
$ws.Range("L113").Value = 2736
$ws.Range("M113").Value = -409.1667000000002
$ws.Range("N113").Value = -7076
$ws.Range("H122").Value = 3609.2727
$ws.Range("I122").Value = 2971.1538
$ws.Range("J122").Value = 4531
$ws.Range("K122").Value = 8913.4614
$ws.Range("L122").Value = 13593
$ws.Range("M122").Value = -6463.4614
$ws.Range("N122").Value = -18493

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4338.3
$ws.Range("I82").Value = 2972.2307
$ws.Range("K82").Value = 2972.2307
$ws.Range("M82").Value = -2611.2307
$ws.Range("H85").Value = 4338.3
$ws.Range("I85").Value = 2972.2307
$ws.Range("K85").Value = 2972.2307
$ws.Range("M85").Value = -1724.2307
$ws.Range("H136").Value = 6079.325
$ws.Range("I136").Value = 5875.3105
$ws.Range("K136").Value = 17625.9315
$ws.Range("M136").Value = -15075.9315

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 271
$ws.Range("J69").Value = 271
$ws.Range("L69").Value = 271
$ws.Range("N69").Value = -1769
$ws.Range("H72").Value = 271
$ws.Range("J72").Value = 271
$ws.Range("L72").Value = 813
$ws.Range("N72").Value = -8301
$ws.Range("H122").Value = 1704.3959
$ws.Range("I122").Value = 1705.7778
$ws.Range("J122").Value = 1700.25
$ws.Range("K122").Value = 5117.3334
$ws.Range("L122").Value = 5100.75
$ws.Range("M122").Value = -2667.3334
$ws.Range("N122").Value = -10000.75
